$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.02
$summary.Range("B4").Value = 0.02
$summary.Range("B5").Value = 0.1
$summary.Range("B6").Value = 4
$summary.Range("B7").Value = 2
$summary.Range("B9").Value = 50

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.02
$status.Range("D4").Value = 4
$status.Range("E4").Value = 0.02
$status.Range("F4").Value = 0.02
$status.Range("G4").Value = 50

# --- New trade row data (used for both "All Trades" and "MarketMaking" sheets) ---
function Add-TradeRow($ws) {
    $ws.Cells.Item(5, 1).Value = 4
    $ws.Cells.Item(5, 2).NumberFormat = "@"
    $ws.Cells.Item(5, 2).Value = "2026-02-17"
    $ws.Cells.Item(5, 3).NumberFormat = "@"
    $ws.Cells.Item(5, 3).Value = "04:05:55"
    $ws.Cells.Item(5, 4).Value = "MarketMaking"
    $ws.Cells.Item(5, 5).Value = "DOWN"
    $ws.Cells.Item(5, 6).Value = 0.78
    $ws.Cells.Item(5, 7).Value = 0.8100000000000001
    $ws.Cells.Item(5, 8).Value = "CLOSED"
    $ws.Cells.Item(5, 9).Value = 3.8462
    $ws.Cells.Item(5, 10).Value = 0.03
    $ws.Cells.Item(5, 11).Value = 100.02
    $ws.Cells.Item(5, 12).Value = 0
    $ws.Cells.Item(5, 13).Value = 0
    $ws.Cells.Item(5, 14).Value = 0.6
    $ws.Cells.Item(5, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(5, 16).Value = "early_exit"
    $ws.Cells.Item(5, 17).Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
